# "novas traducoes e melhorias" - update a couple of job-class translations
# on the Tactics sheet and refresh column C's width to fit the new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tactics")

# hexer row: "Xama" -> "Conjurador Sombrio"
$ws.Range("C9").Value = "Conjurador Sombrio"

# warmage row: "Mago de Guerra" -> "Mago de Batalha"
$ws.Range("C13").Value = "Mago de Batalha"

# Column C now holds longer text, so re-fit its width like Excel would.
$ws.Columns("C").AutoFit()

# Leave the selection where the last edit was made.
$ws.Range("C13").Select()
